$d = $word.ActiveDocument

# 1. Update the ID placeholder text in the first paragraph. Word's Find/Replace
#    merges the two adjacent runs (they share identical formatting), leaving a
#    single run with the new text followed by the original trailing space.
$d.Content.Find.Execute("**ID__AFFARS_5345_topic_2__ID**", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_SUBPART_5345_1__ID**", 2)

# 2. Drop the now-orphaned trailing space that used to be its own run.
$p1 = $d.Paragraphs(1)
$pRange = $p1.Range
$pRange.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark

$trailingSpace = $pRange.Duplicate
$trailingSpace.Start = $pRange.End - 1
$trailingSpace.End = $pRange.End
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}

# 3. Update the paragraph's left indentation (225 twips = 11.25 points).
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# 4. Add a paragraph border (top/left/bottom/right) with 5-twip spacing.
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
